$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 54 (revised quarterly figures) ---
$ws.Range("B54").Value = 106.2
$ws.Range("C54").Value = 89.90000000000001
$ws.Range("D54").Value = 89.5
$ws.Range("E54").Value = 100.3
$ws.Range("F54").Value = 253.1
$ws.Range("G54").Value = 115.5
$ws.Range("H54").Value = 126.9
$ws.Range("I54").Value = 117
$ws.Range("J54").Value = 107.4
$ws.Range("K54").Value = 101.5

# --- Add new row 55 (new quarter: 01-04-2021) ---
# The label "01-04-2021" looks like a date, so a plain Range.Value assignment
# would be auto-converted into a date serial by Excel's type inference. Build
# the text via a TEXT() formula in a scratch cell (forcing a string result),
# then Copy / PasteSpecial the computed value into A55 so the literal text is
# carried across as a shared string without registering a NumberFormat-driven
# cell style (which a direct NumberFormat="@" + Value assignment would leave
# behind, even after clearing it again).
$scratch = $ws.Cells.Item(100, 100)
$scratch.Formula = "=TEXT(DATE(2021,4,1),""dd-mm-yyyy"")"
$scratch.Copy()
$labelCell = $ws.Cells.Item(55, 1)
$labelCell.PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("B55").Value = 103.2
$ws.Range("C55").Value = 99.59999999999999
$ws.Range("D55").Value = 98.90000000000001
$ws.Range("E55").Value = 105.9
$ws.Range("F55").Value = 110.2
$ws.Range("G55").Value = 119.5
$ws.Range("H55").Value = 130.2
$ws.Range("I55").Value = 122
$ws.Range("J55").Value = 106.9
$ws.Range("K55").Value = 106.4
